$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new source/email
$ws.Range("A2").Value = "https://iesab.com.br/preco-do-elevador-residencial/#:~:text=M%C3%A9dia%20de%20pre%C3%A7o%20do%20Elevador,comprimento%2C%20menor%20ser%C3%A1%20o%20custo."
$ws.Range("B2").Value = "emailbit21@gmail.com;"

# Row 3 keeps its site (previously on row 5) but now moves up; set the email value
$ws.Range("A3").Value = "https://coteibem.sindiconet.com.br/fornecedores/manutencao-elevadores/sp/sao-paulo"
$ws.Range("B3").Value = "contato@coteibem.com.br;"

# Remove the now-obsolete rows 4 through 7
$ws.Range("A4:B7").EntireRow.Delete()
